$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The "From Lead" / "To Lead" sample data on row 2 is being appended/updated:
# the old numeric lead ids (10003 / 10005) are replaced with the new demo
# lead names used by the test case. We stage the new text in a scratch cell
# and paste-special (values only) into E2/F2 so the cells keep their
# existing direct formatting (quote-prefixed "text" style) instead of being
# reset to the default style, matching what Excel does when a value is
# retyped in place.
$scratch = $ws.Range("H1")

$scratch.Value = "DemoLeadA"
$scratch.Copy()
$ws.Range("E2").PasteSpecial(-4163)  # xlPasteValues

$scratch.Value = "DemoLeadB"
$scratch.Copy()
$ws.Range("F2").PasteSpecial(-4163)  # xlPasteValues

$scratch.Clear()

# Re-fit the "From Lead"/"To Lead" columns now that their content changed,
# same as Excel does for columns that already use best-fit widths. The two
# columns end up with slightly different best-fit widths (the new values
# are not the same length once rendered), so size them individually.
$ws.Columns.Item(5).AutoFit() | Out-Null
$ws.Columns.Item(6).AutoFit() | Out-Null
$ws.Columns.Item(5).ColumnWidth = 10.666666666666666
$ws.Columns.Item(6).ColumnWidth = 10.833333333333332
